{"js": "// Update the worksheet date and every two-digit x two-digit multiplication\n// prompt in the table to the next day's regenerated set of problems.\nconst replacements = [\n  [\"2025-10-17 Friday\", \"2025-10-18 Saturday\"],\n  [\"55\u00d743=\", \"26\u00d720=\"],\n  [\"41\u00d721=\", \"16\u00d749=\"],\n  [\"71\u00d733=\", \"29\u00d781=\"],\n  [\"66\u00d782=\", \"32\u00d773=\"],\n  [\"76\u00d764=\", \"49\u00d722=\"],\n  [\"92\u00d754=\", \"19\u00d768=\"],\n  [\"78\u00d764=\", \"20\u00d791=\"],\n  [\"54\u00d794=\", \"31\u00d711=\"],\n  [\"27\u00d781=\", \"85\u00d725=\"],\n  [\"73\u00d716=\", \"18\u00d742=\"],\n  [\"86\u00d729=\", \"42\u00d717=\"],\n  [\"53\u00d786=\", \"40\u00d723=\"],\n  [\"59\u00d789=\", \"65\u00d728=\"],\n  [\"74\u00d734=\", \"71\u00d743=\"],\n  [\"20\u00d738=\", \"40\u00d782=\"],\n  [\"26\u00d796=\", \"61\u00d798=\"],\n  [\"33\u00d757=\", \"66\u00d717=\"],\n  [\"67\u00d711=\", \"94\u00d789=\"],\n  [\"78\u00d732=\", \"60\u00d774=\"],\n  [\"32\u00d719=\", \"48\u00d782=\"],\n  [\"21\u00d789=\", \"49\u00d742=\"],\n  [\"99\u00d741=\", \"26\u00d795=\"],\n  [\"34\u00d756=\", \"87\u00d771=\"],\n  [\"71\u00d758=\", \"61\u00d764=\"],\n  [\"46\u00d755=\", \"24\u00d745=\"],\n];\n\nfor (const [findText, newText] of replacements) {\n  const results = context.document.body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n    await context.sync();\n  } else {\n    console.log(`WARNING: text not found -> ${findText}`);\n  }\n}\n", "ps1": "# Update the worksheet date and every two-digit x two-digit multiplication\n# prompt in the table to the next day's regenerated set of problems.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"2025-10-17 Friday\"; Replace = \"2025-10-18 Saturday\" },\n    @{ Find = \"55\u00d743=\";            Replace = \"26\u00d720=\" },\n    @{ Find = \"41\u00d721=\";            Replace = \"16\u00d749=\" },\n    @{ Find = \"71\u00d733=\";            Replace = \"29\u00d781=\" },\n    @{ Find = \"66\u00d782=\";            Replace = \"32\u00d773=\" },\n    @{ Find = \"76\u00d764=\";            Replace = \"49\u00d722=\" },\n    @{ Find = \"92\u00d754=\";            Replace = \"19\u00d768=\" },\n    @{ Find = \"78\u00d764=\";            Replace = \"20\u00d791=\" },\n    @{ Find = \"54\u00d794=\";            Replace = \"31\u00d711=\" },\n    @{ Find = \"27\u00d781=\";            Replace = \"85\u00d725=\" },\n    @{ Find = \"73\u00d716=\";            Replace = \"18\u00d742=\" },\n    @{ Find = \"86\u00d729=\";            Replace = \"42\u00d717=\" },\n    @{ Find = \"53\u00d786=\";            Replace = \"40\u00d723=\" },\n    @{ Find = \"59\u00d789=\";            Replace = \"65\u00d728=\" },\n    @{ Find = \"74\u00d734=\";            Replace = \"71\u00d743=\" },\n    @{ Find = \"20\u00d738=\";            Replace = \"40\u00d782=\" },\n    @{ Find = \"26\u00d796=\";            Replace = \"61\u00d798=\" },\n    @{ Find = \"33\u00d757=\";            Replace = \"66\u00d717=\" },\n    @{ Find = \"67\u00d711=\";            Replace = \"94\u00d789=\" },\n    @{ Find = \"78\u00d732=\";            Replace = \"60\u00d774=\" },\n    @{ Find = \"32\u00d719=\";            Replace = \"48\u00d782=\" },\n    @{ Find = \"21\u00d789=\";            Replace = \"49\u00d742=\" },\n    @{ Find = \"99\u00d741=\";            Replace = \"26\u00d795=\" },\n    @{ Find = \"34\u00d756=\";            Replace = \"87\u00d771=\" },\n    @{ Find = \"71\u00d758=\";            Replace = \"61\u00d764=\" },\n    @{ Find = \"46\u00d755=\";            Replace = \"24\u00d745=\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $pair.Find\n    $range.Find.Replacement.Text = $pair.Replace\n    $found = $range.Find.Execute($pair.Find, $false, $false, $false, $false, $false, $true, 1, $false, $pair.Replace, 2)\n    if (-not $found) {\n        Write-Host \"WARNING: text not found -> $($pair.Find)\"\n    }\n}\n"}
